$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st sheet) - update 想去人数 (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 323
$ws1.Range("F4").Value = 8403
$ws1.Range("F5").Value = 6121
$ws1.Range("F6").Value = 526
$ws1.Range("F7").Value = 110
$ws1.Range("F11").Value = 1090
$ws1.Range("F12").Value = 82

# Sheet "全部类型" (4th sheet) - same underlying rows, different row numbers
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 323
$ws4.Range("F4").Value = 8403
$ws4.Range("F5").Value = 6121
$ws4.Range("F6").Value = 526
$ws4.Range("F7").Value = 110
$ws4.Range("F15").Value = 1090
$ws4.Range("F16").Value = 82
